# Add Address and Email fields to the "users" table row in the Tables sheet,
# and correct its Type / Sync Down Hours / Sync Up Minutes columns to match
# the other "USER"-type (client-synced) table rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tables")

$ws.Range("B24").Value = "name,code,email,phone,address,outstanding_normal,outstanding_overdue,outstanding_critical"
$ws.Range("C24").Value = "USER"
$ws.Range("D24").Value = 1
$ws.Range("E24").Value = 1

# Restore the scroll position / active selection recorded for this sheet.
$win = $excel.ActiveWindow
$win.ScrollRow = 7
$win.ScrollColumn = 1
$ws.Range("B19").Select()
